$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new entry "NIBE Limited" to the next empty row in column A
$ws.Range("A14").Value = "NIBE Limited"

# Update selection to match the new active cell as in the diff
$ws.Range("A14").Select()
